$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the previously-"TODO" cells as done (checkmark) for the
# "Typ + ID" and "Typ + ID + ArgDate(...)" rows.
$check = [char]0x2713
$ws.Range("E7").Value = "$check"
$ws.Range("K7").Value = "$check"
$ws.Range("E8").Value = "$check"
$ws.Range("K8").Value = "$check"

# Move the active selection to J16 (was C18).
$ws.Range("J16").Select()
